# Replace the day's scripture reading content (single paragraph/run) with the
# new reading for 9/22 (Hebrews 1, Isaiah 28, Isaiah 29), per the commit
# "send scriptures to myself".
#
# The whole document is one paragraph containing one run made up of
# alternating <w:t> text nodes and <w:br/> manual line breaks. We rebuild
# that run's text wholesale: each array entry below is one "line" (one
# <w:t> element, or an empty string for a bare <w:br/> with nothing
# between it and its neighbour), joined with the manual-line-break
# character (Chr(11), Word's "^l"), with a trailing break to match the
# document's original trailing <w:br/>.

$lines = @(
    '********************************9月22日读经章节***************************',
    'Chapter 1 of Hebrews',
    '1.神既在古时借着众先知，多次多方的晓谕列祖，',
    '2.就在这末世，借着他儿子晓谕我们，又早已立他为承受万有的，也曾借着他创造诸世界。',
    '3.他是神荣耀所发的光辉，是神本体的真像，常用他权能的命令托住万有，他洗净了人的罪，就坐在高天至大者的右边。',
    '4.他所承受的名，既比天使的名更尊贵，就远超过天使。',
    '5.所有的天使，神从来对那一个说，你是我的儿子，我今日生你。又指着那一个说，我要作他的父，他要作我的子。',
    '6.再者，神使长子到世上来的时候，（或作神再使长子到世上来的时候）就说，神的使者都要拜他。',
    '7.论到使者，又说，神以风为使者，以火焰为仆役。',
    '8.论到子却说，神阿，你的宝座是永永远远的，你的国权是正直的。',
    '9.你喜爱公义，恨恶罪恶。所以神就是你的神，用喜乐油膏你，胜过膏你的同伴。',
    '10.又说，主阿，你起初立了地的根基，天也是你手所造的。',
    '11.天地都要灭没，你却要长存。天地都要像衣服渐渐旧了。',
    '12.你要将天地卷起来，像一件外衣，天地就都改变了。惟有你永不改变，你的年数没有穷尽。',
    '13.所有的天使，神从来对那一个说，你坐在我的右边，等我使你仇敌作你的脚凳。',
    '14.天使岂不都是服役的灵，奉差遣为那将要承受救恩的人效力吗？',
    'Chapter 28 of Isaiah',
    '1.祸哉，以法莲的酒徒，住在肥美谷的山上，他们心里高傲，以所夸的为冠冕，犹如将残之花。',
    '2.看哪，主有一大能大力者，像一阵冰雹，像毁灭的暴风，像涨溢的大水，他必用手将冠冕摔落于地。',
    '3.以法莲高傲的酒徒，他的冠冕，必被踏在脚下。',
    '4.那荣美将残之花，就是在肥美谷山上的，必像夏令以前初熟的无花果。看见这果的就注意，一到手中就吞吃了。',
    '5.到那日，万军之耶和华必作他余剩之民的荣冠华冕。',
    '6.也作了在位上行审判者公平之灵，并城门口打退仇敌者的力量。',
    '7.就是这地的人，也因酒摇摇晃晃，因浓酒东倒西歪。祭司和先知因浓酒摇摇晃晃，被酒所困，因浓酒东倒西歪。他们错解默示，谬行审判。',
    '8.因为各席上满了呕吐的污秽，无一处干净。',
    '9.讥诮先知的说，他要将知识指教谁呢？要使谁明白传言呢？是那刚断奶离怀的吗？',
    '10.他竟命上加命，令上加令，律上加律，例上加例，这里一点，那里一点。',
    '11.先知说，不然，主要借异邦人的嘴唇，和外邦人的舌头，对这百姓说话。',
    '12.他曾对他们说，你们要使疲乏人得安息。这样才得安息，才得舒畅。他们却不肯听。',
    '13.所以耶和华向他们说的话，是命上加命，令上加令，律上加律，例上加例，这里一点，那里一点，以致他们前行仰面跌倒，而且跌碎，并陷入网罗，被缠住。',
    '14.所以你们这些亵慢的人，就是辖管住在耶路撒冷这百姓的，要听耶和华的话。',
    '15.你们曾说，我们与死亡立约，与阴间结盟。敌军（原文作鞭子）如水涨漫经过的时候，必不临到我们。因我们以谎言为避所，在虚假以下藏身。',
    '16.所以主耶和华如此说，看哪，我在锡安放一块石头，作为根基，是试验过的石头，是稳固根基，宝贵的房角石，信靠的人必不着急。',
    '17.我必以公平为准绳，以公义为线铊。冰雹必冲去谎言的避所，大水必漫过藏身之处。',
    '18.你们与死亡所立的约，必然废掉，与阴间所结的盟，必立不住。敌军（原文作鞭子）如水涨漫经过的时候，你们必被他践踏。',
    '19.每逢经过必将你们掳去。因为每早晨他必经过，白昼黑夜都必如此。明白传言的，必受惊恐。',
    '20.原来床榻短，使人不能舒身。被窝窄，使人不能遮体。',
    '21.耶和华必兴起，像在毗拉心山，他必发怒，像在基遍谷，好作成他的工，就是非常的工，成就他的事，就是奇异的事。',
    '22.现在你们不可亵慢，恐怕捆你们的绑索更结实了。因为我从主万军之耶和华那里听见已经决定，在全地上施行灭绝的事。',
    '23.你们当侧耳听我的声音，留心听我的言语。',
    '24.那耕地为要撒种的，岂是常常耕地呢？岂是常常开垦耙地呢？',
    '25.他拉平了地面，岂不就撒种小茴香，播种大茴香，按行列种小麦，在定处种大麦，在田边种粗麦呢？',
    '26.因为他的神教导他务农相宜，并且指教他。',
    '27.原来打小茴香，不用尖利的器具，轧大茴香，也不用碌碡。（原文作车轮下同）但用杖打小茴香，用棍打大茴香。',
    '28.作饼的粮食是用磨磨碎，因它不必常打。虽用碌碡和马打散，却不磨它。',
    '29.这也是出于万军之耶和华。他的谋略奇妙，他的智慧广大。',
    'Chapter 29 of Isaiah',
    '1.唉，亚利伊勒，亚利伊勒，大卫安营的城。任凭你年上加年，节期照常周流。',
    '2.我终必使亚利伊勒困难。他必悲伤哀号，我却仍以他为亚利伊勒。',
    '3.我必四围安营攻击你，屯兵围困你，筑垒攻击你。',
    '4.你必败落，从地中说话。你的言语必微细出于尘埃。你的声音必像那交鬼者的声音出于地。你的言语低低微微出于尘埃。',
    '5.你仇敌的群众，却要像细尘，强暴人的群众，也要像飞糠。这事必顷刻之间忽然临到。',
    '6.万军之耶和华必用雷轰，地震，大声，旋风，暴风，并吞灭的火焰，向他讨罪。',
    '7.那时，攻击亚利伊勒列国的群众，就是一切攻击亚利伊勒和他的保障，并使他困难的，必如梦景，如夜间的异象。',
    '8.又必像饥饿的人，梦中吃饭，醒了仍觉腹空。或像口渴的人，梦中喝水，醒了仍觉发昏，心里想喝。攻击锡安山列国的群众，也必如此。',
    '9.你们等候惊奇吧。你们宴乐昏迷吧。他们醉了，却非因酒。他们东倒西歪，却非因浓酒。',
    '10.因为耶和华将沉睡的灵，浇灌你们，封闭你们的眼，蒙盖你们的头。你们的眼，就是先知。你们的头，就是先见。',
    '11.所有的默示你们看如封住的书卷。人将这书卷交给识字的，说，请念吧。他说，我不能念，因为是封住了。',
    '12.又将这书卷交给不识字的人，说，请念吧。他说，我不识字。',
    '13.主说，因为这百姓亲近我，用嘴唇尊敬我，心却远离我。他们敬畏我，不过是领受人的吩咐。',
    '14.所以我在这百姓中要行奇妙的事，就是奇妙又奇妙的事。他们智慧人的智慧，必然消灭，聪明人的聪明，必然隐藏。',
    '15.祸哉，那些向耶和华深藏谋略的，又在暗中行事，说，谁看见我们呢？谁知道我们呢？',
    '16.你们把事颠倒了，岂可看窑匠如泥吗？被制作的物，岂可论制作物的说，他没有制作我。或是被创造的物论造物的说，他没有聪明？',
    '17.利巴嫩变为肥田，肥田看如树林，不是只有一点点时候吗？',
    '18.那时，聋子必听见这书上的话。瞎子的眼，必从迷蒙黑暗中得以看见。',
    '19.谦卑人，必因耶和华增添欢喜。人间贫穷的，必因以色列的圣者快乐。',
    '20.因为强暴人已归无有。亵慢人已经灭绝。一切找机会作孽的，都被剪除。',
    '21.他们在争讼的事上，定无罪的为有罪，为城门口责备人的，设下网罗，用虚无的事，屈枉义人。',
    '',
    '22.所以救赎亚伯拉罕的耶和华，论雅各家如此说，雅各必不再羞愧，面容也不至变色。',
    '23.但他看见他的众子，就是我手的工作，在他那里，他们必尊我的名为圣，必尊雅各的圣者为圣，必敬畏以色列的神。',
    '24.心中迷糊的，必得明白，发怨言的，必受训诲。'
)

$d = $word.ActiveDocument
$d.Content.Text = ($lines -join [char]11) + [char]11
